$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label in B1 (was "inn_real", now "A091RC1Q027SBEA")
$ws.Range("B1").Value = "A091RC1Q027SBEA"

# Rescale B2:B51 values by dividing by 10,000,000 (convert from raw units to millions/other scale)
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Value2
    $cell.Value = $old / 10000000
}
